$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 5.197675333333334
$ws.Range("N2").Value = 15.593026
$ws.Range("O2").Value = 0.4667706154782232
$ws.Range("P2").Value = 0.4667706154782231
$ws.Range("Q2").Value = 0.7484028758960001
$ws.Range("R2").Value = 6.735625883064001
$ws.Range("S2").Value = 0.4667706154782232
$ws.Range("T2").Value = 0.4667706154782231

$ws.Range("O3").Value = 0.1721501820052131
$ws.Range("P3").Value = 0.172150182005213
$ws.Range("S3").Value = 0.1721501820052131
$ws.Range("T3").Value = 0.172150182005213

$ws.Range("M4").Value = 0.9113383333333335
$ws.Range("N4").Value = 2.734015
$ws.Range("O4").Value = 0.08184157868246321
$ws.Range("P4").Value = 0.08184157868246318
$ws.Range("Q4").Value = 0.13122178394
$ws.Range("R4").Value = 1.18099605546
$ws.Range("S4").Value = 0.08184157868246321
$ws.Range("T4").Value = 0.08184157868246318

$ws.Range("M5").Value = 0.9066646666666666
$ws.Range("N5").Value = 2.719994
$ws.Range("O5").Value = 0.08142186599811185
$ws.Range("P5").Value = 0.08142186599811183
$ws.Range("Q5").Value = 0.130548832024
$ws.Range("R5").Value = 1.174939488216
$ws.Range("S5").Value = 0.08142186599811185
$ws.Range("T5").Value = 0.08142186599811183

$ws.Range("M6").Value = 1.443486333333333
$ws.Range("N6").Value = 4.330459
$ws.Range("O6").Value = 0.1296304522761144
$ws.Range("P6").Value = 0.1296304522761144
$ws.Range("Q6").Value = 0.207844710164
$ws.Range("R6").Value = 1.870602391476
$ws.Range("S6").Value = 0.1296304522761144
$ws.Range("T6").Value = 0.1296304522761144

$ws.Range("M7").Value = 0.7592703333333333
$ws.Range("N7").Value = 2.277811
$ws.Range("O7").Value = 0.06818530555987445
$ws.Range("P7").Value = 0.06818530555987443
$ws.Range("Q7").Value = 0.109325816756
$ws.Range("R7").Value = 0.983932350804
$ws.Range("S7").Value = 0.06818530555987445
$ws.Range("T7").Value = 0.06818530555987443
